$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: split the run(s) covering [start,end) away from their neighbours by
# toggling a character-formatting property on and back off. Word (and this
# interop layer) always materialises a fresh <w:r> boundary around a range
# whose formatting is written to, even when the net value doesn't change,
# which lets us recreate the original run layout after a text edit has
# coalesced runs together.
# ---------------------------------------------------------------------------
function Split-At($start, $end) {
    $r = $d.Range($start, $end)
    $r.Font.Bold = $true
    $r.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# Helper: force a real (non no-op) text replacement over [start,end) so that
# the interop layer actually rewrites/merges the underlying runs, even when
# the desired final text is identical to the concatenation of the existing
# runs' text (a same-text assignment is treated as a no-op and left alone).
# ---------------------------------------------------------------------------
function Retype-Range($start, $end, $finalText) {
    $tmp = "@@@@@@@@@@"
    $r = $d.Range($start, $end)
    $r.Text = $tmp
    $r2 = $d.Range($start, $start + $tmp.Length)
    $r2.Text = $finalText
}

# ===========================================================================
# 1. "Last updated 9.00am, 2 June 2019." -> "Last updated 8.50am, 9 June 2019."
# ===========================================================================

$timeRange = $d.Content
$timeRange.Find.Execute("9.00")
$timeStart = $timeRange.Start
$timeEnd = $timeRange.End
$timeRange.Text = "8.50"

$spaceRange = $d.Content
$spaceRange.Find.Execute(", ")
$spaceRange.Find.Execute("2 June")
$dateWordRange = $d.Content
$dateWordRange.Find.Execute("2 June")
$dateStart = $dateWordRange.Start
$dateEnd = $dateWordRange.End

$spaceBeforeRange = $d.Range($dateStart - 1, $dateStart)
$spaceStart = $spaceBeforeRange.Start
$spaceEnd = $spaceBeforeRange.End

# " " -> " 9" (the run right before "2 June")
$spaceBeforeRange.Text = " 9"

# "2 June" -> " June" (shifted by the +1 character just inserted above)
$dateRange = $d.Range($dateStart + 1, $dateEnd + 1)
$dateRange.Text = " June"

# Recompute the run boundaries of the whole sentence so it matches the
# original, un-merged run layout (only the "8.50", " 9" and " June" runs'
# text should actually differ from before).
$sentence = $d.Content
$sentence.Find.Execute("Last updated")
$lastUpdatedStart = $sentence.Start
$lastUpdatedEnd = $sentence.End

$periodSearch = $d.Content
$periodSearch.Find.Execute("2019.")
$periodEnd = $periodSearch.End

Split-At $lastUpdatedEnd ($lastUpdatedEnd + 1)          # " " before "8.50"
Split-At ($lastUpdatedEnd + 1) ($lastUpdatedEnd + 5)    # "8.50"
Split-At ($lastUpdatedEnd + 5) ($lastUpdatedEnd + 7)    # "am"
Split-At ($lastUpdatedEnd + 7) ($lastUpdatedEnd + 8)    # ","
Split-At ($lastUpdatedEnd + 8) ($lastUpdatedEnd + 10)   # " 9"
Split-At ($lastUpdatedEnd + 10) ($lastUpdatedEnd + 15)  # " June"
Split-At ($lastUpdatedEnd + 15) ($lastUpdatedEnd + 16)  # " " before "2019"
Split-At ($lastUpdatedEnd + 16) ($periodEnd - 1)        # "2019"
Split-At ($periodEnd - 1) $periodEnd                    # "."

# Move the (hidden) _GoBack bookmark so it sits right after the new " 9" run,
# i.e. right before " June" -- matching where Word recorded the last edit.
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$bookmarkPos = $d.Range(($lastUpdatedEnd + 10), ($lastUpdatedEnd + 10))
$d.Bookmarks.Add("_GoBack", $bookmarkPos)

# ===========================================================================
# 2. "100 Carols for Cho" + bookmark + "irs" -> single run "100 Carols for
#    Choirs" (the bookmark used to live here; it moved to the date above).
# ===========================================================================

$choirsRange = $d.Content
$choirsRange.Find.Execute("100 Carols for Choirs")
$choirsStart = $choirsRange.Start
$choirsEnd = $choirsRange.End

Retype-Range $choirsStart $choirsEnd "100 Carols for Choirs"

Write-Host "done"
